# Applies scheduled-runner market/profit data updates across all sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1067
$ws.Range("I11").Value = 1067
$ws.Range("K11").Value = 1067
$ws.Range("M11").Value = -927
$ws.Range("H15").Value = 490.825
$ws.Range("I15").Value = 490.825
$ws.Range("K15").Value = 1472.475
$ws.Range("M15").Value = -1303.475
$ws.Range("H17").Value = 2048
$ws.Range("J17").Value = 3030
$ws.Range("L17").Value = 9090
$ws.Range("N17").Value = -9426
$ws.Range("H33").Value = 609.9474
$ws.Range("I33").Value = 178.58333
$ws.Range("J33").Value = 1349.4286
$ws.Range("K33").Value = 178.58333
$ws.Range("L33").Value = 1349.4286
$ws.Range("M33").Value = 50.41667000000001
$ws.Range("N33").Value = -1807.4286
$ws.Range("H43").Value = 6750.5
$ws.Range("I43").Value = 5001
$ws.Range("K43").Value = 5001
$ws.Range("M43").Value = -4932
$ws.Range("H95").Value = 23500
$ws.Range("J95").Value = 23500
$ws.Range("L95").Value = 23500
$ws.Range("N95").Value = -28992
$ws.Range("H111").Value = 2375.8
$ws.Range("I111").Value = 1693
$ws.Range("K111").Value = 5079
$ws.Range("M111").Value = -2012
$ws.Range("H112").Value = 3838.3076
$ws.Range("I112").Value = 1199
$ws.Range("K112").Value = 3597
$ws.Range("M112").Value = -2489
$ws.Range("H125").Value = 3240.25
$ws.Range("I125").Value = 1050
$ws.Range("K125").Value = 9450
$ws.Range("M125").Value = -6990
$ws.Range("H127").Value = 834
$ws.Range("I127").Value = 834
$ws.Range("K127").Value = 2502
$ws.Range("M127").Value = 2458
$ws.Range("H131").Value = 2000
$ws.Range("J131").Value = 2000
$ws.Range("L131").Value = 6000
$ws.Range("N131").Value = -16080
$ws.Range("H137").Value = 1884.8
$ws.Range("I137").Value = 1474.6666
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 4423.9998
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -1873.9998
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 5276.5
$ws.Range("I138").Value = 4222.3076
$ws.Range("K138").Value = 12666.9228
$ws.Range("M138").Value = -7526.9228

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1577.6364
$ws.Range("I122").Value = 1535.4
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4606.200000000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2156.200000000001
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 3039
$ws.Range("I132").Value = 2648.75
$ws.Range("K132").Value = 7946.25
$ws.Range("M132").Value = -5416.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7500.5
$ws.Range("I86").Value = 7500.5
$ws.Range("K86").Value = 7500.5
$ws.Range("M86").Value = -6377.5
$ws.Range("H89").Value = 7500.5
$ws.Range("I89").Value = 7500.5
$ws.Range("K89").Value = 37502.5
$ws.Range("M89").Value = -31886.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6385.1
$ws.Range("I4").Value = 1284.3334
$ws.Range("J4").Value = 8571.143
$ws.Range("K4").Value = 1284.3334
$ws.Range("L4").Value = 8571.143
$ws.Range("M4").Value = -1172.3334
$ws.Range("N4").Value = -8795.143
$ws.Range("H43").Value = 30583.166
$ws.Range("J43").Value = 30583.166
$ws.Range("L43").Value = 30583.166
$ws.Range("N43").Value = -30951.166
$ws.Range("H58").Value = 5861.25
$ws.Range("I58").Value = 5035.6665
$ws.Range("K58").Value = 5035.6665
$ws.Range("M58").Value = -4832.6665
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H101").Value = 30583.166
$ws.Range("J101").Value = 30583.166
$ws.Range("L101").Value = 30583.166
$ws.Range("N101").Value = -37073.166
$ws.Range("H107").Value = 564.05884
$ws.Range("I107").Value = 199.08333
$ws.Range("K107").Value = 199.08333
$ws.Range("M107").Value = 1720.91667
$ws.Range("H132").Value = 6078.4
$ws.Range("I132").Value = 3873.75
$ws.Range("K132").Value = 11621.25
$ws.Range("M132").Value = -9091.25
$ws.Range("H136").Value = 5861.25
$ws.Range("I136").Value = 5035.6665
$ws.Range("K136").Value = 15106.9995
$ws.Range("M136").Value = -12556.9995
$ws.Range("H138").Value = 1249.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 302.66666
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H55").Value = 891.46155
$ws.Range("I55").Value = 399.75
$ws.Range("J55").Value = 1110
$ws.Range("K55").Value = 1199.25
$ws.Range("L55").Value = 3330
$ws.Range("M55").Value = -1022.25
$ws.Range("N55").Value = -3684
$ws.Range("H98").Value = 568.4286
$ws.Range("I98").Value = 526
$ws.Range("K98").Value = 1578
$ws.Range("M98").Value = -80
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H124").Value = 15000
$ws.Range("J124").Value = 15000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -54820
$ws.Range("H135").Value = 302.66666
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H139").Value = 2242.6
$ws.Range("J139").Value = 6000
$ws.Range("L139").Value = 18000
$ws.Range("N139").Value = -28280
$ws.Range("H140").Value = 1193.75
$ws.Range("I140").Value = 925
$ws.Range("K140").Value = 2775
$ws.Range("M140").Value = 2405

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 918.4
$ws.Range("I102").Value = 918.4
$ws.Range("K102").Value = 918.4
$ws.Range("M102").Value = 703.6
$ws.Range("H126").Value = 7279.8
$ws.Range("I126").Value = 5750
$ws.Range("J126").Value = 8299.666999999999
$ws.Range("K126").Value = 17250
$ws.Range("L126").Value = 24899.001
$ws.Range("M126").Value = -14780
$ws.Range("N126").Value = -29839.001
$ws.Range("H132").Value = 3727
$ws.Range("I132").Value = 3682.6
$ws.Range("K132").Value = 11047.8
$ws.Range("M132").Value = -8517.799999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 33333.168
$ws.Range("I74").Value = 26666.666
$ws.Range("K74").Value = 26666.666
$ws.Range("M74").Value = -25668.666
$ws.Range("H77").Value = 33333.168
$ws.Range("I77").Value = 26666.666
$ws.Range("K77").Value = 79999.99800000001
$ws.Range("M77").Value = -75007.99800000001
$ws.Range("H136").Value = 4852.1816
$ws.Range("I136").Value = 4210.5713
$ws.Range("K136").Value = 12631.7139
$ws.Range("M136").Value = -10081.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20000
$ws.Range("I2").Value = 20000
$ws.Range("K2").Value = 20000
$ws.Range("M2").Value = -19888
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 700
$ws.Range("K122").Value = 2100
$ws.Range("M122").Value = 350
$ws.Range("H132").Value = 1990.2307
$ws.Range("I132").Value = 1897.5454
$ws.Range("K132").Value = 5692.6362
$ws.Range("M132").Value = -3162.6362
$ws.Range("H136").Value = 3297.739
$ws.Range("I136").Value = 2579.353
$ws.Range("K136").Value = 7738.059
$ws.Range("M136").Value = -5188.059
